$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MTSI")

# Insert two new columns before column D, shifting existing D:K data to F:M.
$colD = $ws.Range("D1")
$colD.EntireColumn.Insert()
$colD.EntireColumn.Insert()

# New values for the newly-inserted columns D (most recent quarter) and E (prior quarter).
$updates = @(
    @{Row=7;   D=43462;   E=43371}
    @{Row=8;   D=150700;  E=151200}
    @{Row=9;   D=74100;   E=80200}
    @{Row=10;  D=76600;   E=71000}
    @{Row=12;  D=43500;   E=46200}
    @{Row=13;  D=0;       E=0}
    @{Row=14;  D=5000;    E=3500}
    @{Row=15;  D=12500;   E=12400}
    @{Row=17;  D=165100;  E=168700}
    @{Row=18;  D=-14400;  E=-17500}
    @{Row=20;  D=900;     E=-900}
    @{Row=21;  D=14700;   E=10300}
    @{Row=22;  D=8800;    E=8100}
    @{Row=23;  D=-22300;  E=-26400}
    @{Row=24;  D=1100;    E=-9200}
    @{Row=25;  D=0;       E=0}
    @{Row=26;  D=-23400;  E=-17200}
    @{Row=27;  D=-28900;  E=7700}
    @{Row=28;  D=0;       E=0}
    @{Row=29;  D=0;       E=700}
    @{Row=30;  D=0;       E=0}
    @{Row=31;  D=0;       E=0}
    @{Row=32;  D=-900;    E=900}
    @{Row=33;  D=-28900;  E=8400}
    @{Row=34;  D=0;       E=0}
    @{Row=35;  D=-28900;  E=8400}
    @{Row=38;  D=43462;   E=43371}
    @{Row=41;  D=87100;   E=94700}
    @{Row=42;  D=98700;   E=98200}
    @{Row=43;  D=127500;  E=115000}
    @{Row=44;  D=120900;  E=122800}
    @{Row=45;  D=25300;   E=28200}
    @{Row=46;  D=459400;  E=458900}
    @{Row=47;  D=26500;   E=31100}
    @{Row=48;  D=148900;  E=149900}
    @{Row=49;  D=806700;  E=826900}
    @{Row=50;  D=0;       E=0}
    @{Row=51;  D=0;       E=0}
    @{Row=52;  D=16300;   E=15800}
    @{Row=53;  D=0;       E=0}
    @{Row=54;  D=1457800; E=1482500}
    @{Row=57;  D=32600;   E=42000}
    @{Row=58;  D=8000;    E=7400}
    @{Row=59;  D=52600;   E=57700}
    @{Row=60;  D=93200;   E=107000}
    @{Row=61;  D=688900;  E=687400}
    @{Row=62;  D=18500;   E=19400}
    @{Row=63;  D=0;       E=0}
    @{Row=64;  D=0;       E=0}
    @{Row=65;  D=0;       E=0}
    @{Row=66;  D=800600;  E=813800}
    @{Row=68;  D=0;       E=0}
    @{Row=69;  D=0;       E=0}
    @{Row=70;  D=0;       E=0}
    @{Row=71;  D=0;       E=0}
    @{Row=72;  D=-431400; E=-408000}
    @{Row=73;  D=0;       E=0}
    @{Row=74;  D=0;       E=0}
    @{Row=75;  D=0;       E=0}
    @{Row=76;  D=657300;  E=668700}
    @{Row=77;  D=0;       E=0}
    @{Row=80;  D=43462;   E=43371}
    @{Row=81;  D=-28900;  E=8400}
    @{Row=83;  D=28200;   E=28700}
    @{Row=84;  D=0;       E=0}
    @{Row=85;  D=0;       E=0}
    @{Row=86;  D=0;       E=0}
    @{Row=87;  D=0;       E=0}
    @{Row=88;  D=0;       E=0}
    @{Row=89;  D=2900;    E=25100}
    @{Row=91;  D=-11500;  E=-13600}
    @{Row=92;  D=0;       E=0}
    @{Row=93;  D=0;       E=0}
    @{Row=94;  D=-12300;  E=-13500}
    @{Row=96;  D=0;       E=0}
    @{Row=97;  D=0;       E=0}
    @{Row=98;  D=0;       E=0}
    @{Row=99;  D=0;       E=0}
    @{Row=100; D=1900;    E=-2000}
    @{Row=101; D=0;       E=-200}
    @{Row=102; D=-7600;   E=9400}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
